$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 8-10 (matrix rows for MuSCs->* self pairs previously present)
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Update remaining data rows (2-7) with the new TPM-derived values

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga11"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.419591
$ws.Range("H2").Value = 4.258773
$ws.Range("I2").Value = 0.001848767113890483
$ws.Range("J2").Value = 0.001848767113890483
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 33.94639966666666
$ws.Range("N2").Value = 101.839199
$ws.Range("O2").Value = 0.9746097333921855
$ws.Range("P2").Value = 0.9746097333921855
$ws.Range("Q2").Value = 48.19000344920299
$ws.Range("R2").Value = 433.7100310428269
$ws.Range("S2").Value = 0.001801826423973044
$ws.Range("T2").Value = 0.001801826423973044

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga11"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.419591
$ws.Range("H3").Value = 4.258773
$ws.Range("I3").Value = 0.001848767113890483
$ws.Range("J3").Value = 0.001848767113890483
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8843623333333334
$ws.Range("N3").Value = 2.653087
$ws.Range("O3").Value = 0.02539026660781448
$ws.Range("P3").Value = 0.02539026660781448
$ws.Range("Q3").Value = 1.255432809139
$ws.Range("R3").Value = 11.298895282251
$ws.Range("S3").Value = [double]"4.694068991743908E-05"
$ws.Range("T3").Value = [double]"4.694068991743908E-05"

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga11"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 735.4993083333334
$ws.Range("H4").Value = 2206.497925
$ws.Range("I4").Value = 0.9578582377148513
$ws.Range("J4").Value = 0.9578582377148513
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 33.94639966666666
$ws.Range("N4").Value = 101.839199
$ws.Range("O4").Value = 0.9746097333921855
$ws.Range("P4").Value = 0.9746097333921855
$ws.Range("Q4").Value = 24967.55347524023
$ws.Range("R4").Value = 224707.9812771621
$ws.Range("S4").Value = 0.9335379616867799
$ws.Range("T4").Value = 0.9335379616867799

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga11"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 735.4993083333334
$ws.Range("H5").Value = 2206.497925
$ws.Range("I5").Value = 0.9578582377148513
$ws.Range("J5").Value = 0.9578582377148513
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8843623333333334
$ws.Range("N5").Value = 2.653087
$ws.Range("O5").Value = 0.02539026660781448
$ws.Range("P5").Value = 0.02539026660781448
$ws.Range("Q5").Value = 650.4478844827195
$ws.Range("R5").Value = 5854.030960344476
$ws.Range("S5").Value = 0.02432027602807141
$ws.Range("T5").Value = 0.02432027602807141

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga11"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 30.939307
$ws.Range("H6").Value = 92.81792100000001
$ws.Range("I6").Value = 0.04029299517125823
$ws.Range("J6").Value = 0.04029299517125823
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 33.94639966666666
$ws.Range("N6").Value = 101.839199
$ws.Range("O6").Value = 0.9746097333921855
$ws.Range("P6").Value = 0.9746097333921855
$ws.Range("Q6").Value = 1050.278080831698
$ws.Range("R6").Value = 9452.50272748528
$ws.Range("S6").Value = 0.0392699452814326
$ws.Range("T6").Value = 0.0392699452814326

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga11"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 30.939307
$ws.Range("H7").Value = 92.81792100000001
$ws.Range("I7").Value = 0.04029299517125823
$ws.Range("J7").Value = 0.04029299517125823
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8843623333333334
$ws.Range("N7").Value = 2.653087
$ws.Range("O7").Value = 0.02539026660781448
$ws.Range("P7").Value = 0.02539026660781448
$ws.Range("Q7").Value = 27.36155773023634
$ws.Range("R7").Value = 246.254019572127
$ws.Range("S7").Value = 0.001023049889825628
$ws.Range("T7").Value = 0.001023049889825628
